$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 454 (before the current row 455),
# shifting all the existing data (old rows 455-482) down to rows 457-484.
$ws.Rows("455:456").Insert()

# New row 455: "Primera" record for 2023-04-25 (serial 45041)
$ws.Cells.Item(455, 1).Value = 11
$ws.Cells.Item(455, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(455, 3).Value = "Bíobío"
$ws.Cells.Item(455, 4).Value = 45041
$ws.Cells.Item(455, 5).Value = 8
$ws.Cells.Item(455, 6).Value = 100112023
$ws.Cells.Item(455, 7).Value = "Brócoli"
$ws.Cells.Item(455, 8).Value = "Sin especificar"
$ws.Cells.Item(455, 9).Value = "Primera"
$ws.Cells.Item(455, 10).Value = 2000
$ws.Cells.Item(455, 11).Value = 700
$ws.Cells.Item(455, 12).Value = 800
$ws.Cells.Item(455, 13).Value = 750
$ws.Cells.Item(455, 14).Value = "$/unidad"
$ws.Cells.Item(455, 15).Value = "Región Metropolitana"
$ws.Cells.Item(455, 16).Value = 750
$ws.Cells.Item(455, 17).Value = 1
$ws.Cells.Item(455, 18).Value = "Hortaliza"

# New row 456: "Segunda" record for 2023-04-25 (serial 45041)
$ws.Cells.Item(456, 1).Value = 11
$ws.Cells.Item(456, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(456, 3).Value = "Bíobío"
$ws.Cells.Item(456, 4).Value = 45041
$ws.Cells.Item(456, 5).Value = 8
$ws.Cells.Item(456, 6).Value = 100112023
$ws.Cells.Item(456, 7).Value = "Brócoli"
$ws.Cells.Item(456, 8).Value = "Sin especificar"
$ws.Cells.Item(456, 9).Value = "Segunda"
$ws.Cells.Item(456, 10).Value = 1500
$ws.Cells.Item(456, 11).Value = 600
$ws.Cells.Item(456, 12).Value = 600
$ws.Cells.Item(456, 13).Value = 600
$ws.Cells.Item(456, 14).Value = "$/unidad"
$ws.Cells.Item(456, 15).Value = "Región Metropolitana"
$ws.Cells.Item(456, 16).Value = 600
$ws.Cells.Item(456, 17).Value = 1
$ws.Cells.Item(456, 18).Value = "Hortaliza"
